# Scheduled-runner update: refresh cached market-price-derived values
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/
# CUL/GSM/LTW/WVR leve-profit tables, matching the latest Universalis pull.
# Each cell was literal (no formulas), so values are rewritten in place; a
# few rows gain/lose their LeveProfitNQ/HQ (M/N) cell because profit flips sign
# relative to the break-even it's computed against.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 914.2857
$ws.Cells.Item(38, 8).Value = 632.0909
$ws.Cells.Item(38, 9).Value = 300.75
$ws.Cells.Item(38, 10).Value = 821.4286
$ws.Cells.Item(38, 11).Value = 902.25
$ws.Cells.Item(38, 12).Value = 2464.2858
$ws.Cells.Item(38, 13).Value = -530.25
$ws.Cells.Item(38, 14).Value = -3208.2858
$ws.Cells.Item(128, 8).Value = 20780
$ws.Cells.Item(128, 10).Value = 20780
$ws.Cells.Item(128, 12).Value = 20780
$ws.Cells.Item(128, 14).Value = -30740
$ws.Cells.Item(129, 8).Value = 860
$ws.Cells.Item(129, 10).Value = 860
$ws.Cells.Item(129, 12).Value = 2580
$ws.Cells.Item(129, 14).Value = -12580
$ws.Cells.Item(135, 8).Value = 23818052
$ws.Cells.Item(135, 9).Value = 807
$ws.Cells.Item(135, 10).Value = 100033240
$ws.Cells.Item(135, 11).Value = 7263
$ws.Cells.Item(135, 12).Value = 900299160
$ws.Cells.Item(135, 13).Value = -4728
$ws.Cells.Item(135, 14).Value = -900304230
$ws.Cells.Item(137, 8).Value = 2448.3
$ws.Cells.Item(137, 9).Value = 2276.6
$ws.Cells.Item(137, 10).Value = 2620
$ws.Cells.Item(137, 11).Value = 6829.799999999999
$ws.Cells.Item(137, 12).Value = 7860
$ws.Cells.Item(137, 13).Value = -4279.799999999999
$ws.Cells.Item(137, 14).Value = -12960
$ws.Cells.Item(138, 8).Value = 1936.9773
$ws.Cells.Item(138, 9).Value = 553.2727
$ws.Cells.Item(138, 10).Value = 3320.682
$ws.Cells.Item(138, 11).Value = 1659.8181
$ws.Cells.Item(138, 12).Value = 9962.045999999998
$ws.Cells.Item(138, 13).Value = 3480.1819
$ws.Cells.Item(138, 14).Value = -20242.046
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2669.2
$ws.Cells.Item(2, 10).Value = 4000.6667
$ws.Cells.Item(2, 12).Value = 4000.6667
$ws.Cells.Item(2, 14).Value = -4226.6667
$ws.Cells.Item(4, 8).Value = 63.333332
$ws.Cells.Item(4, 10).Value = 45
$ws.Cells.Item(4, 12).Value = 45
$ws.Cells.Item(4, 14).Value = -277
$ws.Cells.Item(116, 8).Value = 2669.2
$ws.Cells.Item(116, 10).Value = 4000.6667
$ws.Cells.Item(116, 12).Value = 4000.6667
$ws.Cells.Item(116, 14).Value = -8588.6667
$ws.Cells.Item(132, 8).Value = 11409.02
$ws.Cells.Item(132, 9).Value = 1158.9556
$ws.Cells.Item(132, 11).Value = 3476.8668
$ws.Cells.Item(132, 13).Value = -946.8667999999998
$ws.Cells.Item(139, 8).Value = 50640
$ws.Cells.Item(139, 10).Value = 50640
$ws.Cells.Item(139, 12).Value = 50640
$ws.Cells.Item(139, 14).Value = -60920
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2669.2
$ws.Cells.Item(3, 10).Value = 4000.6667
$ws.Cells.Item(3, 12).Value = 4000.6667
$ws.Cells.Item(3, 14).Value = -4228.6667
$ws.Cells.Item(20, 8).Value = 1855.2354
$ws.Cells.Item(20, 9).Value = 2103.25
$ws.Cells.Item(20, 11).Value = 2103.25
$ws.Cells.Item(20, 13).Value = -1856.25
$ws.Cells.Item(22, 8).Value = 282.14285
$ws.Cells.Item(22, 9).Value = 282.14285
$ws.Cells.Item(22, 11).Value = 282.14285
$ws.Cells.Item(22, 13).Value = -109.14285
$ws.Cells.Item(107, 8).Value = 392.85715
$ws.Cells.Item(107, 9).Value = 392.85715
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 392.85715
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 1527.14285
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 43000
$ws.Cells.Item(132, 10).Value = 43000
$ws.Cells.Item(132, 12).Value = 43000
$ws.Cells.Item(132, 14).Value = -53120
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(62, 8).Value = 6602.4
$ws.Cells.Item(62, 9).Value = 5000
$ws.Cells.Item(62, 10).Value = 7003
$ws.Cells.Item(62, 11).Value = 5000
$ws.Cells.Item(62, 12).Value = 7003
$ws.Cells.Item(62, 13).Value = -4376
$ws.Cells.Item(62, 14).Value = -8251
$ws.Cells.Item(65, 8).Value = 6602.4
$ws.Cells.Item(65, 9).Value = 5000
$ws.Cells.Item(65, 10).Value = 7003
$ws.Cells.Item(65, 11).Value = 25000
$ws.Cells.Item(65, 12).Value = 35015
$ws.Cells.Item(65, 13).Value = -21880
$ws.Cells.Item(65, 14).Value = -41255
$ws.Cells.Item(94, 8).Value = 2933
$ws.Cells.Item(94, 9).Value = 2148.2222
$ws.Cells.Item(94, 10).Value = 3575.0908
$ws.Cells.Item(94, 11).Value = 2148.2222
$ws.Cells.Item(94, 12).Value = 3575.0908
$ws.Cells.Item(94, 13).Value = -1697.2222
$ws.Cells.Item(94, 14).Value = -4477.0908
$ws.Cells.Item(99, 8).Value = 16207120
$ws.Cells.Item(99, 9).Value = 3207646.5
$ws.Cells.Item(99, 11).Value = 3207646.5
$ws.Cells.Item(99, 13).Value = -3206148.5
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(126, 8).Value = 16207120
$ws.Cells.Item(126, 9).Value = 3207646.5
$ws.Cells.Item(126, 11).Value = 9622939.5
$ws.Cells.Item(126, 13).Value = -9620469.5
$ws.Cells.Item(132, 8).Value = 2827.8215
$ws.Cells.Item(132, 9).Value = 1889.2609
$ws.Cells.Item(132, 11).Value = 5667.7827
$ws.Cells.Item(132, 13).Value = -3137.7827
$ws.Cells.Item(134, 8).Value = 1136.2727
$ws.Cells.Item(134, 9).Value = 791.4286
$ws.Cells.Item(134, 10).Value = 1739.75
$ws.Cells.Item(134, 11).Value = 2374.2858
$ws.Cells.Item(134, 12).Value = 5219.25
$ws.Cells.Item(134, 13).Value = 160.7142000000003
$ws.Cells.Item(134, 14).Value = -10289.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 283
$ws.Cells.Item(33, 9).Value = 249.5
$ws.Cells.Item(33, 10).Value = 299.75
$ws.Cells.Item(33, 11).Value = 1497
$ws.Cells.Item(33, 12).Value = 1798.5
$ws.Cells.Item(33, 13).Value = -1214
$ws.Cells.Item(33, 14).Value = -2364.5
$ws.Cells.Item(74, 8).Value = 8300
$ws.Cells.Item(74, 10).Value = 9850
$ws.Cells.Item(74, 12).Value = 29550
$ws.Cells.Item(74, 14).Value = -31672
$ws.Cells.Item(77, 8).Value = 8300
$ws.Cells.Item(77, 10).Value = 9850
$ws.Cells.Item(77, 12).Value = 88650
$ws.Cells.Item(77, 14).Value = -99258
$ws.Cells.Item(92, 8).Value = 900
$ws.Cells.Item(92, 9).Value = 533.3333
$ws.Cells.Item(92, 11).Value = 1599.9999
$ws.Cells.Item(92, 13).Value = -351.9999
$ws.Cells.Item(113, 8).Value = 10518.546
$ws.Cells.Item(113, 9).Value = 100000
$ws.Cells.Item(113, 10).Value = 1570.4
$ws.Cells.Item(113, 11).Value = 300000
$ws.Cells.Item(113, 12).Value = 4711.200000000001
$ws.Cells.Item(113, 13).Value = -297830
$ws.Cells.Item(113, 14).Value = -9051.200000000001
$ws.Cells.Item(131, 8).Value = 776.75
$ws.Cells.Item(131, 10).Value = 791.8261
$ws.Cells.Item(131, 12).Value = 2375.4783
$ws.Cells.Item(131, 14).Value = -12455.4783
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3912.5293
$ws.Cells.Item(80, 10).Value = 4088
$ws.Cells.Item(80, 12).Value = 4088
$ws.Cells.Item(80, 14).Value = -6084
$ws.Cells.Item(83, 8).Value = 3912.5293
$ws.Cells.Item(83, 10).Value = 4088
$ws.Cells.Item(83, 12).Value = 20440
$ws.Cells.Item(83, 14).Value = -30424
$ws.Cells.Item(132, 8).Value = 26000.545
$ws.Cells.Item(132, 9).Value = 3669.3845
$ws.Cells.Item(132, 10).Value = 58256.668
$ws.Cells.Item(132, 11).Value = 11008.1535
$ws.Cells.Item(132, 12).Value = 174770.004
$ws.Cells.Item(132, 13).Value = -8478.1535
$ws.Cells.Item(132, 14).Value = -179830.004
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3612.75
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 14).Value = -2090
$ws.Cells.Item(27, 8).Value = 3612.75
$ws.Cells.Item(27, 10).Value = 1500
$ws.Cells.Item(27, 12).Value = 1500
$ws.Cells.Item(27, 14).Value = -1714
$ws.Cells.Item(61, 8).Value = 3631.25
$ws.Cells.Item(61, 9).Value = 1507.5
$ws.Cells.Item(61, 10).Value = 10002.5
$ws.Cells.Item(61, 11).Value = 1507.5
$ws.Cells.Item(61, 12).Value = 10002.5
$ws.Cells.Item(61, 13).Value = -1305.5
$ws.Cells.Item(61, 14).Value = -10406.5
$ws.Cells.Item(93, 8).Value = 3563.25
$ws.Cells.Item(93, 9).Value = 3505.9
$ws.Cells.Item(93, 10).Value = 3850
$ws.Cells.Item(93, 11).Value = 3505.9
$ws.Cells.Item(93, 12).Value = 3850
$ws.Cells.Item(93, 13).Value = -2257.9
$ws.Cells.Item(93, 14).Value = -6346
$ws.Cells.Item(113, 8).Value = 3631.25
$ws.Cells.Item(113, 9).Value = 1507.5
$ws.Cells.Item(113, 10).Value = 10002.5
$ws.Cells.Item(113, 11).Value = 1507.5
$ws.Cells.Item(113, 12).Value = 10002.5
$ws.Cells.Item(113, 13).Value = 662.5
$ws.Cells.Item(113, 14).Value = -14342.5
$ws.Cells.Item(122, 8).Value = 1786426
$ws.Cells.Item(122, 9).Value = 2453873.5
$ws.Cells.Item(122, 10).Value = 6566.3335
$ws.Cells.Item(122, 11).Value = 7361620.5
$ws.Cells.Item(122, 12).Value = 19699.0005
$ws.Cells.Item(122, 13).Value = -7359170.5
$ws.Cells.Item(122, 14).Value = -24599.0005
$ws.Cells.Item(132, 8).Value = 2499
$ws.Cells.Item(132, 9).Value = 1783.4286
$ws.Cells.Item(132, 11).Value = 5350.2858
$ws.Cells.Item(132, 13).Value = -2820.2858
$ws.Cells.Item(136, 8).Value = 2139.1667
$ws.Cells.Item(136, 9).Value = 1766.6666
$ws.Cells.Item(136, 10).Value = 2325.4167
$ws.Cells.Item(136, 11).Value = 5299.9998
$ws.Cells.Item(136, 12).Value = 6976.250100000001
$ws.Cells.Item(136, 13).Value = -2749.9998
$ws.Cells.Item(136, 14).Value = -12076.2501
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 12990
$ws.Cells.Item(33, 9).Value = 12990
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 12990
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -12740
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(36, 8).Value = 12990
$ws.Cells.Item(36, 9).Value = 12990
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 12990
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -12740
$ws.Cells.Item(36, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 18869180
$ws.Cells.Item(136, 9).Value = 27778906
$ws.Cells.Item(136, 10).Value = 1526.0588
$ws.Cells.Item(136, 11).Value = 83336718
$ws.Cells.Item(136, 12).Value = 4578.1764
$ws.Cells.Item(136, 13).Value = -83334168
$ws.Cells.Item(136, 14).Value = -9678.1764
